$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Luciano -> Marília Lucí, new phone, new date (01/07/2025)
$ws.Range("A2").Value = "Marília Lucí"
$ws.Range("B2").Value = 558388008228
$ws.Range("C2").Value = 45839

# Row 3 (Ismar pai) stays the same - no change needed

# Row 4: Rejania keeps her name, phone number updated
$ws.Range("B4").Value = 5583987592215

# Row 5: Augusto removed - clear the whole row's contents
$ws.Range("A5:C5").ClearContents()

# Update the active selection to D4
$ws.Range("D4").Select()
